$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Add($ws.Range("A6"), "https://leetcode.com/problems/product-of-array-except-self/", "", "", "https://leetcode.com/problems/product-of-array-except-self/") | Out-Null
$ws.Range("A6").Value = "Product of Array Except Self"
$ws.Range("B6").Value = "Product of any element will be product of all elements to the right of it and to the left of it (2 times iterarte over the array one from beginning one from end)"

$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B:B").ColumnWidth = 139.45

$ws.Range("B15").Select()
